# Display_Bookshelves/BookshelvesDetails.xlsx edit:
#  1. Insert two new rows (new rows 2 and 3) for a new "Rhodes Folding Book Shelf"
#     product (Teak and Mahogany variants), pushing the existing product rows down.
#  2. Re-order the "Boeberg Cabinet Inserts" (Dark Walnut Finish) row so it sits
#     directly after the "Boeberg Cabinet Inserts" (Columbian Walnut Finish) row,
#     ahead of "Boeberg Drawer Inserts".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two rows right below the header row for the new product ---
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "Rhodes Folding Book Shelf"
$ws.Range("B2").Value = "(Teak Finish, Tall Configuration, 60 Book Book Capacity)"
$ws.Range("C2").Value = "₹15,599"

$ws.Range("A3").Value = "Rhodes Folding Book Shelf"
$ws.Range("B3").Value = "(Mahogany Finish, Tall Configuration, 60 Book Book Capacity)"
$ws.Range("C3").Value = "₹15,599"

# --- 2. Swap the "Boeberg Drawer Inserts" row and the "(Dark Walnut Finish)" row ---
# After the insert above, the table (rows 7-9) now reads:
#   row 7: Boeberg Cabinet Inserts | (Columbian Walnut Finish) | ₹1,699 ₹1,019
#   row 8: Boeberg Drawer Inserts  | (Columbian Walnut Finish) | ₹2,699 ₹1,619
#   row 9: Boeberg Cabinet Inserts | (Dark Walnut Finish)      | ₹1,499 ₹899
# Target order swaps rows 8 and 9 so the two "Boeberg Cabinet Inserts" rows are
# adjacent, followed by "Boeberg Drawer Inserts".

$ws.Range("A8").Value = "Boeberg Cabinet Inserts"
$ws.Range("B8").Value = "(Dark Walnut Finish)"
$ws.Range("C8").Value = "₹1,499 ₹899"

$ws.Range("A9").Value = "Boeberg Drawer Inserts"
$ws.Range("B9").Value = "(Columbian Walnut Finish)"
$ws.Range("C9").Value = "₹2,699 ₹1,619"
